$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.618.36'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '2.343.71'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.524'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.96%  '
$ws.Range('D9').Value = '2.342.76'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.28'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('D15').Value = '2.768.16'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = '60.540.41'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '2.343.50'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '317.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  -3.69%  '
$ws.Range('E26').Value = '  +9.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '496.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('E29').Value = '  -3.60%  '
$ws.Range('D30').Value = '0.0₃0853'
$ws.Range('E30').Value = '  -6.47%  '
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').Value = '  -2.23%  '
$ws.Range('E33').Value = '  -2.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.375'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.96%  '
$ws.Range('E37').Value = '  +3.35%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.11%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.19'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.38%  '
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '141.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('E47').Value = '  -5.17%  '
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0898'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('E50').Value = '  -1.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.97%  '
